$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the merged-AHB header labels: "<Feld>_old" -> "<Feld>_FV2210",
#    "<Feld>_new" -> "<Feld>_FV2304" (row 1 header cells only).
$headerRange = $ws.Range("A1:U1")
for ($c = 1; $c -le $headerRange.Columns.Count; $c++) {
    $cell = $headerRange.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2210"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2304"
        }
    }
}

# 2) Turn the data range into a proper Excel table ("Table1") with an
#    autofilter, matching the regenerated merged-AHB export.
$dataRange = $ws.Range("A1:U60")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# 3) Freeze the header row so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
